# Weekly update: a new price record for "Femacal de La Calera - Berenjena"
# is inserted as row 227 (pushing the former rows 227:239 down to 228:240).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 227, shifting existing rows 227-239 -> 228-240
$ws.Rows("227:227").Insert()

# Populate the newly inserted row with the new record
$ws.Cells.Item(227, 1).Value = 3
$ws.Cells.Item(227, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(227, 3).Value = "Coquimbo"
$ws.Cells.Item(227, 4).Value = 44610
$ws.Cells.Item(227, 5).Value = 5
$ws.Cells.Item(227, 6).Value = 100112001
$ws.Cells.Item(227, 7).Value = "Berenjena"
$ws.Cells.Item(227, 8).Value = "Sin especificar"
$ws.Cells.Item(227, 9).Value = "Primera"
$ws.Cells.Item(227, 10).Value = 110
$ws.Cells.Item(227, 11).Value = 9500
$ws.Cells.Item(227, 12).Value = 10000
$ws.Cells.Item(227, 13).Value = 9727
$ws.Cells.Item(227, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(227, 15).Value = "Región Metropolitana"
$ws.Cells.Item(227, 16).Value = 162
$ws.Cells.Item(227, 17).Value = 60
$ws.Cells.Item(227, 18).Value = "Hortaliza"
